$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 295, pushing the existing rows 295-349 down to 296-350.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new data point.
$ws.Range("A295").Value = 3
$ws.Range("B295").Value = "Femacal de La Calera"
$ws.Range("C295").Value = "Coquimbo"
$ws.Range("D295").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D295").Value = 44694
$ws.Range("E295").Value = 5
$ws.Range("F295").Value = 100114013
$ws.Range("G295").Value = "Zanahoria"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 340
$ws.Range("K295").Value = 7000
$ws.Range("L295").Value = 7500
$ws.Range("M295").Value = 7235
$ws.Range("N295").Value = "`$/saco 20 kilos"
$ws.Range("O295").Value = "Provincia de Quillota"
$ws.Range("P295").Value = 362
$ws.Range("Q295").Value = 20
$ws.Range("R295").Value = "Hortaliza"
